# "sl added in duplicate rows"
#
# The "Stock Report" sheet stored several damage rows as a single
# serialized-looking string in column B (an inline string that looked like
# a Ruby array dump: `[nil, nil, ..., "FLOORS-(F)", "Threshold plate-(...)",
# "description", "", "", ""]`). This script explodes that payload into the
# real columns (C..AB), mirroring the layout already used by the sibling
# "real" rows, and clears column B back out. It also fixes up column widths
# that were sized for the old one-giant-cell layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# Column B was sized to hold the giant serialized string; now that the data
# is split across C:AB it goes back to a normal (still hidden) narrow width.
$ws.Columns.Item(2).ColumnWidth = 7.857142857142857
$ws.Columns.Item(2).Hidden = $true

# Columns X (24) and Y (25) grow to fit the longer part-name / description
# text that now lives in real cells instead of the packed string.
$ws.Columns.Item(24).ColumnWidth = 46.42857142857143
$ws.Columns.Item(25).ColumnWidth = 58.42857142857143

# Row number -> [Damage Area Name, Damage Part Name, Damage Description]
# (columns Z/AA/AB - Damage Component/Type/Repair Type - are blank, same as
# the already-normal rows such as row 7/8/10/12/15).
$rows = @{
    9  = @("FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY BADLY SODA DUST & ODOUR .")
    11 = @("FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY BADLY SODA DUST & ODOUR .")
    13 = @("DOORS-(D)", "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)", 'RIGHT DOOR MIDDLE GASKET CUT 3"')
    14 = @("PANELS-(PANELS)", "Right side panel -(Right side panel )", 'LEFT SIDE PANEL 3,4 CUT 2" & 4"X2"')
    16 = @("PANELS-(PANELS)", "Right side panel -(Right side panel )", "EXTERNAL PANEL RUSTED CORROSION")
    17 = @("FLOORS-(F)", "FLOOR BOARD-(FLOOR BOARD)", "FLOOR BOARD DIRTY BY WOOD & MUD DUST & SCRATCHED")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    # Drop the packed/serialized payload out of column B.
    $ws.Cells.Item($r, 2).ClearContents()

    # Materialise empty cells C..V (3..22), matching the normal rows'
    # layout, so the row looks like every other data row structurally.
    for ($c = 3; $c -le 22; $c++) {
        $ws.Cells.Item($r, $c).Font.Name = "Arial"
    }

    # M (13) and U (21) - Gate In Date / Di Date - carry the date format
    # used by the rest of the sheet even though they stay blank here.
    $ws.Cells.Item($r, 13).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 21).NumberFormat = "YYYY-MM-DD"

    # W, X, Y - Damage Area Name / Damage Part Name / Damage Description.
    $ws.Cells.Item($r, 23).Value = $vals[0]
    $ws.Cells.Item($r, 24).Value = $vals[1]
    $ws.Cells.Item($r, 25).Value = $vals[2]

    # Z, AA, AB - Damage Component / Damage Type / Repair Type - blank.
    $ws.Cells.Item($r, 26).Font.Name = "Arial"
    $ws.Cells.Item($r, 27).Font.Name = "Arial"
    $ws.Cells.Item($r, 28).Font.Name = "Arial"
}
